# Generate Report for Handoff
# Rename the localized e2e markdown file (and its generated xliff handoff/handback
# artifacts) from the 1b57126d... id to the 179a4f7d... id, and bump the
# associated timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldId = "1b57126d-5c64-47e0-8545-d9bb9a6d7efc"
$newId = "179a4f7d-981c-4bd1-bd6b-91f36f1ec11e"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

$oldMdPath = "e2e\$oldId.md"
$newMdPath = "e2e\$newId.md"

$oldZhXlf = "$oldId.5345b12ede8d953177f294180262dec7c981c22a.zh-cn.xlf"
$newZhXlf = "$newId.7dfa765555422fb89b36af70038e7aab27e5c8f7.zh-cn.xlf"

$oldDeXlf = "$oldId.5345b12ede8d953177f294180262dec7c981c22a.de-de.xlf"
$newDeXlf = "$newId.7dfa765555422fb89b36af70038e7aab27e5c8f7.de-de.xlf"

$oldGenDate = "2016-08-18 08:59:10"
$newGenDate = "2016-08-18 08:59:34"

$oldZhDate = "2016-08-18 08:58:59"
$newZhDate = "2016-08-18 08:59:29"

# The hyperlink target address is unaffected by this change (same commit/URL
# in the upstream repo); only the displayed text changes. Since the engine
# does not surface the stored Address of hyperlinks loaded from the original
# file, hardcode the existing target so Delete+Add keeps it unchanged.
$hlAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28af2fff4d7dd0a707efdcf4f67f54c1d2da3b34/e2e/$oldId.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newGenDate

# Update the hyperlink display text on B2 while preserving its target address.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hlAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newMdPath)

# ---- zh-cn sheet ----
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhDate

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hlAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newMd)

# ---- de-de sheet ----
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newGenDate

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hlAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newMd)
